$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A107: correct the timestamp (R script re-run produced a revised value)
$ws.Range("A107").Value = 45476.2916666667

# Append new row 108 with the latest data pulled by the R script
# Carry the date-time number format from the cell above by copying formats first.
$ws.Range("A107").Copy()
$ws.Range("A108").PasteSpecial(-4122)

$ws.Range("A108").Value = 45477.6404398148
$ws.Range("B108").Value = 50038
$ws.Range("C108").Value = 0.720000028610229
$ws.Range("D108").Value = 0.644999980926514
$ws.Range("E108").Value = 0.675000011920929
$ws.Range("F108").Value = 0.665000021457672

# G (adj_close) is stored as TEXT in this sheet (matches the rest of the
# column), even though its content looks numeric. Build it via a formula
# result (a real string) in a scratch cell, then paste-values it in so the
# destination cell picks up text type instead of being reinterpreted as a
# number - plain `.Value = "0.665000021457672"` would otherwise coerce to
# a numeric cell.
$ws.Range("Z1").Formula = "=""0.665000021457672"""
$ws.Range("Z1").Copy()
$ws.Range("G108").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("H108").Value = "BWZ.MI"
